$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2236597029243569
$ws.Range("C2").Value = 0.6761979329963813
$ws.Range("D2").Value = 0.6208552414790156
$ws.Range("E2").Value = 0.7879436791287913
$ws.Range("F2").Value = 0.7762373633456771
$ws.Range("G2").Value = 19

$ws.Range("B3").Value = 0.1243738929906032
$ws.Range("C3").Value = 0.5670298149555084
$ws.Range("D3").Value = 0.4487500652448925
$ws.Range("E3").Value = 0.6698880990470666
$ws.Range("F3").Value = 0.677324390080757
$ws.Range("G3").Value = 18

$ws.Range("B4").Value = 0.1692601144216072
$ws.Range("C4").Value = 0.4938470918437637
$ws.Range("D4").Value = 0.3483999257463576
$ws.Range("E4").Value = 0.5902541196352276
$ws.Range("F4").Value = 0.5828682296416691
$ws.Range("G4").Value = 17

$ws.Range("B5").Value = 0.2998499526236316
$ws.Range("C5").Value = 0.5009693588263471
$ws.Range("D5").Value = 0.3362117046188728
$ws.Range("E5").Value = 0.5798376536746065
$ws.Range("F5").Value = 0.5125639711936881
$ws.Range("G5").Value = 16

$ws.Range("B6").Value = 0.3249303377857756
$ws.Range("C6").Value = 0.510714471671429
$ws.Range("D6").Value = 0.3463237248587889
$ws.Range("E6").Value = 0.5884927568447966
$ws.Range("F6").Value = 0.5078779385610399
$ws.Range("G6").Value = 15

$ws.Range("B7").Value = 0.357871698984652
$ws.Range("C7").Value = 0.5111964517621793
$ws.Range("D7").Value = 0.3510459376640385
$ws.Range("E7").Value = 0.592491297542874
$ws.Range("F7").Value = 0.4900261363687482
$ws.Range("G7").Value = 14

$ws.Range("B8").Value = 0.3633620203708059
$ws.Range("C8").Value = 0.5275006984044196
$ws.Range("D8").Value = 0.3717216992278948
$ws.Range("E8").Value = 0.6096898385473509
$ws.Range("F8").Value = 0.5095722583647349
$ws.Range("G8").Value = 13

$ws.Range("B9").Value = 0.4177687993451564
$ws.Range("C9").Value = 0.536905318813658
$ws.Range("D9").Value = 0.3885450245204918
$ws.Range("E9").Value = 0.623333798634802
$ws.Range("F9").Value = 0.4831874337779735
$ws.Range("G9").Value = 12

$ws.Range("B10").Value = 0.435803605847665
$ws.Range("C10").Value = 0.5614457896225986
$ws.Range("D10").Value = 0.409650469768552
$ws.Range("E10").Value = 0.6400394282921577
$ws.Range("F10").Value = 0.4916281680178602
$ws.Range("G10").Value = 11

$ws.Range("B11").Value = 0.4121559002167278
$ws.Range("C11").Value = 0.5537197283926004
$ws.Range("D11").Value = 0.4082939408543487
$ws.Range("E11").Value = 0.638978826608792
$ws.Range("F11").Value = 0.5146967335462779
$ws.Range("G11").Value = 10
